$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-06-21 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-06-22 Thursday", 2) | Out-Null
$d.Content.Find.Execute("34×88=", $true, $false, $false, $false, $false, $true, 1, $false, "97×90=", 2) | Out-Null
$d.Content.Find.Execute("37×42=", $true, $false, $false, $false, $false, $true, 1, $false, "42×35=", 2) | Out-Null
$d.Content.Find.Execute("38×19=", $true, $false, $false, $false, $false, $true, 1, $false, "94×24=", 2) | Out-Null
$d.Content.Find.Execute("87×43=", $true, $false, $false, $false, $false, $true, 1, $false, "76×30=", 2) | Out-Null
$d.Content.Find.Execute("32×73=", $true, $false, $false, $false, $false, $true, 1, $false, "62×94=", 2) | Out-Null
$d.Content.Find.Execute("42×40=", $true, $false, $false, $false, $false, $true, 1, $false, "72×10=", 2) | Out-Null
$d.Content.Find.Execute("19×90=", $true, $false, $false, $false, $false, $true, 1, $false, "67×31=", 2) | Out-Null
$d.Content.Find.Execute("65×65=", $true, $false, $false, $false, $false, $true, 1, $false, "72×59=", 2) | Out-Null
$d.Content.Find.Execute("25×22=", $true, $false, $false, $false, $false, $true, 1, $false, "36×91=", 2) | Out-Null
$d.Content.Find.Execute("81×95=", $true, $false, $false, $false, $false, $true, 1, $false, "49×56=", 2) | Out-Null
$d.Content.Find.Execute("55×25=", $true, $false, $false, $false, $false, $true, 1, $false, "88×37=", 2) | Out-Null
$d.Content.Find.Execute("36×45=", $true, $false, $false, $false, $false, $true, 1, $false, "28×51=", 2) | Out-Null
$d.Content.Find.Execute("74×84=", $true, $false, $false, $false, $false, $true, 1, $false, "51×51=", 2) | Out-Null
$d.Content.Find.Execute("30×69=", $true, $false, $false, $false, $false, $true, 1, $false, "61×50=", 2) | Out-Null
$d.Content.Find.Execute("24×72=", $true, $false, $false, $false, $false, $true, 1, $false, "94×52=", 2) | Out-Null
$d.Content.Find.Execute("14×16=", $true, $false, $false, $false, $false, $true, 1, $false, "29×10=", 2) | Out-Null
$d.Content.Find.Execute("19×98=", $true, $false, $false, $false, $false, $true, 1, $false, "32×24=", 2) | Out-Null
$d.Content.Find.Execute("47×50=", $true, $false, $false, $false, $false, $true, 1, $false, "34×41=", 2) | Out-Null
$d.Content.Find.Execute("63×58=", $true, $false, $false, $false, $false, $true, 1, $false, "25×45=", 2) | Out-Null
$d.Content.Find.Execute("41×93=", $true, $false, $false, $false, $false, $true, 1, $false, "16×53=", 2) | Out-Null
$d.Content.Find.Execute("88×78=", $true, $false, $false, $false, $false, $true, 1, $false, "58×15=", 2) | Out-Null
$d.Content.Find.Execute("93×83=", $true, $false, $false, $false, $false, $true, 1, $false, "67×13=", 2) | Out-Null
$d.Content.Find.Execute("36×49=", $true, $false, $false, $false, $false, $true, 1, $false, "59×40=", 2) | Out-Null
$d.Content.Find.Execute("18×100=", $true, $false, $false, $false, $false, $true, 1, $false, "79×35=", 2) | Out-Null
$d.Content.Find.Execute("62×91=", $true, $false, $false, $false, $false, $true, 1, $false, "56×82=", 2) | Out-Null
$d.Content.Find.Execute("93×82=", $true, $false, $false, $false, $false, $true, 1, $false, "27×65=", 2) | Out-Null
$d.Content.Find.Execute("98×22=", $true, $false, $false, $false, $false, $true, 1, $false, "74×13=", 2) | Out-Null
$d.Content.Find.Execute("71×28=", $true, $false, $false, $false, $false, $true, 1, $false, "10×87=", 2) | Out-Null
$d.Content.Find.Execute("11×76=", $true, $false, $false, $false, $false, $true, 1, $false, "69×46=", 2) | Out-Null
$d.Content.Find.Execute("14×11=", $true, $false, $false, $false, $false, $true, 1, $false, "97×81=", 2) | Out-Null
$d.Content.Find.Execute("53×98=", $true, $false, $false, $false, $false, $true, 1, $false, "84×34=", 2) | Out-Null
$d.Content.Find.Execute("20×92=", $true, $false, $false, $false, $false, $true, 1, $false, "16×64=", 2) | Out-Null
$d.Content.Find.Execute("78×80=", $true, $false, $false, $false, $false, $true, 1, $false, "79×40=", 2) | Out-Null
$d.Content.Find.Execute("53×91=", $true, $false, $false, $false, $false, $true, 1, $false, "55×52=", 2) | Out-Null
$d.Content.Find.Execute("68×20=", $true, $false, $false, $false, $false, $true, 1, $false, "44×100=", 2) | Out-Null
$d.Content.Find.Execute("100×88=", $true, $false, $false, $false, $false, $true, 1, $false, "76×36=", 2) | Out-Null
$d.Content.Find.Execute("39×14=", $true, $false, $false, $false, $false, $true, 1, $false, "40×61=", 2) | Out-Null
$d.Content.Find.Execute("90×70=", $true, $false, $false, $false, $false, $true, 1, $false, "35×57=", 2) | Out-Null
$d.Content.Find.Execute("26×76=", $true, $false, $false, $false, $false, $true, 1, $false, "55×93=", 2) | Out-Null
$d.Content.Find.Execute("98×41=", $true, $false, $false, $false, $false, $true, 1, $false, "30×57=", 2) | Out-Null
$d.Content.Find.Execute("70×78=", $true, $false, $false, $false, $false, $true, 1, $false, "71×70=", 2) | Out-Null
$d.Content.Find.Execute("98×61=", $true, $false, $false, $false, $false, $true, 1, $false, "32×20=", 2) | Out-Null
$d.Content.Find.Execute("100×76=", $true, $false, $false, $false, $false, $true, 1, $false, "62×15=", 2) | Out-Null
$d.Content.Find.Execute("33×59=", $true, $false, $false, $false, $false, $true, 1, $false, "85×45=", 2) | Out-Null
$d.Content.Find.Execute("53×66=", $true, $false, $false, $false, $false, $true, 1, $false, "76×32=", 2) | Out-Null
$d.Content.Find.Execute("89×15=", $true, $false, $false, $false, $false, $true, 1, $false, "35×28=", 2) | Out-Null
$d.Content.Find.Execute("84×48=", $true, $false, $false, $false, $false, $true, 1, $false, "100×67=", 2) | Out-Null
$d.Content.Find.Execute("14×55=", $true, $false, $false, $false, $false, $true, 1, $false, "82×94=", 2) | Out-Null
$d.Content.Find.Execute("28×79=", $true, $false, $false, $false, $false, $true, 1, $false, "42×49=", 2) | Out-Null
$d.Content.Find.Execute("61×17=", $true, $false, $false, $false, $false, $true, 1, $false, "87×98=", 2) | Out-Null
$d.Content.Find.Execute("55×59=", $true, $false, $false, $false, $false, $true, 1, $false, "88×31=", 2) | Out-Null
$d.Content.Find.Execute("49×51=", $true, $false, $false, $false, $false, $true, 1, $false, "30×19=", 2) | Out-Null
$d.Content.Find.Execute("67×38=", $true, $false, $false, $false, $false, $true, 1, $false, "11×66=", 2) | Out-Null
$d.Content.Find.Execute("65×67=", $true, $false, $false, $false, $false, $true, 1, $false, "36×28=", 2) | Out-Null
$d.Content.Find.Execute("58×98=", $true, $false, $false, $false, $false, $true, 1, $false, "74×70=", 2) | Out-Null
$d.Content.Find.Execute("28×20=", $true, $false, $false, $false, $false, $true, 1, $false, "96×14=", 2) | Out-Null
$d.Content.Find.Execute("76×72=", $true, $false, $false, $false, $false, $true, 1, $false, "44×30=", 2) | Out-Null
$d.Content.Find.Execute("40×90=", $true, $false, $false, $false, $false, $true, 1, $false, "22×94=", 2) | Out-Null
$d.Content.Find.Execute("21×57=", $true, $false, $false, $false, $false, $true, 1, $false, "67×78=", 2) | Out-Null
$d.Content.Find.Execute("100×50=", $true, $false, $false, $false, $false, $true, 1, $false, "87×41=", 2) | Out-Null
$d.Content.Find.Execute("58×24=", $true, $false, $false, $false, $false, $true, 1, $false, "71×50=", 2) | Out-Null
$d.Content.Find.Execute("94×50=", $true, $false, $false, $false, $false, $true, 1, $false, "22×93=", 2) | Out-Null
$d.Content.Find.Execute("33×88=", $true, $false, $false, $false, $false, $true, 1, $false, "60×81=", 2) | Out-Null
$d.Content.Find.Execute("61×65=", $true, $false, $false, $false, $false, $true, 1, $false, "53×55=", 2) | Out-Null
$d.Content.Find.Execute("92×28=", $true, $false, $false, $false, $false, $true, 1, $false, "61×36=", 2) | Out-Null
$d.Content.Find.Execute("80×79=", $true, $false, $false, $false, $false, $true, 1, $false, "22×86=", 2) | Out-Null
$d.Content.Find.Execute("33×77=", $true, $false, $false, $false, $false, $true, 1, $false, "61×35=", 2) | Out-Null
$d.Content.Find.Execute("69×88=", $true, $false, $false, $false, $false, $true, 1, $false, "67×14=", 2) | Out-Null
$d.Content.Find.Execute("13×80=", $true, $false, $false, $false, $false, $true, 1, $false, "88×54=", 2) | Out-Null
$d.Content.Find.Execute("43×80=", $true, $false, $false, $false, $false, $true, 1, $false, "76×51=", 2) | Out-Null
$d.Content.Find.Execute("83×28=", $true, $false, $false, $false, $false, $true, 1, $false, "62×85=", 2) | Out-Null
$d.Content.Find.Execute("20×67=", $true, $false, $false, $false, $false, $true, 1, $false, "10×15=", 2) | Out-Null
$d.Content.Find.Execute("67×69=", $true, $false, $false, $false, $false, $true, 1, $false, "21×22=", 2) | Out-Null
$d.Content.Find.Execute("49×66=", $true, $false, $false, $false, $false, $true, 1, $false, "57×55=", 2) | Out-Null
$d.Content.Find.Execute("66×65=", $true, $false, $false, $false, $false, $true, 1, $false, "64×80=", 2) | Out-Null
$d.Content.Find.Execute("77×29=", $true, $false, $false, $false, $false, $true, 1, $false, "72×31=", 2) | Out-Null
$d.Content.Find.Execute("18×23=", $true, $false, $false, $false, $false, $true, 1, $false, "60×92=", 2) | Out-Null
$d.Content.Find.Execute("57×19=", $true, $false, $false, $false, $false, $true, 1, $false, "52×59=", 2) | Out-Null
$d.Content.Find.Execute("99×39=", $true, $false, $false, $false, $false, $true, 1, $false, "91×66=", 2) | Out-Null
$d.Content.Find.Execute("51×49=", $true, $false, $false, $false, $false, $true, 1, $false, "35×34=", 2) | Out-Null
$d.Content.Find.Execute("61×99=", $true, $false, $false, $false, $false, $true, 1, $false, "36×63=", 2) | Out-Null
$d.Content.Find.Execute("17×93=", $true, $false, $false, $false, $false, $true, 1, $false, "56×48=", 2) | Out-Null
$d.Content.Find.Execute("27×41=", $true, $false, $false, $false, $false, $true, 1, $false, "20×11=", 2) | Out-Null
$d.Content.Find.Execute("83×11=", $true, $false, $false, $false, $false, $true, 1, $false, "38×94=", 2) | Out-Null
$d.Content.Find.Execute("88×41=", $true, $false, $false, $false, $false, $true, 1, $false, "72×13=", 2) | Out-Null
$d.Content.Find.Execute("38×99=", $true, $false, $false, $false, $false, $true, 1, $false, "23×93=", 2) | Out-Null
$d.Content.Find.Execute("85×64=", $true, $false, $false, $false, $false, $true, 1, $false, "85×83=", 2) | Out-Null
$d.Content.Find.Execute("67×58=", $true, $false, $false, $false, $false, $true, 1, $false, "63×48=", 2) | Out-Null
$d.Content.Find.Execute("22×47=", $true, $false, $false, $false, $false, $true, 1, $false, "83×76=", 2) | Out-Null
$d.Content.Find.Execute("24×35=", $true, $false, $false, $false, $false, $true, 1, $false, "62×72=", 2) | Out-Null
$d.Content.Find.Execute("18×19=", $true, $false, $false, $false, $false, $true, 1, $false, "51×95=", 2) | Out-Null
$d.Content.Find.Execute("47×69=", $true, $false, $false, $false, $false, $true, 1, $false, "79×78=", 2) | Out-Null
$d.Content.Find.Execute("78×85=", $true, $false, $false, $false, $false, $true, 1, $false, "80×11=", 2) | Out-Null
$d.Content.Find.Execute("68×74=", $true, $false, $false, $false, $false, $true, 1, $false, "76×18=", 2) | Out-Null
$d.Content.Find.Execute("22×34=", $true, $false, $false, $false, $false, $true, 1, $false, "68×17=", 2) | Out-Null
$d.Content.Find.Execute("98×67=", $true, $false, $false, $false, $false, $true, 1, $false, "89×81=", 2) | Out-Null
$d.Content.Find.Execute("72×33=", $true, $false, $false, $false, $false, $true, 1, $false, "57×88=", 2) | Out-Null
$d.Content.Find.Execute("11×27=", $true, $false, $false, $false, $false, $true, 1, $false, "44×85=", 2) | Out-Null
$d.Content.Find.Execute("59×61=", $true, $false, $false, $false, $false, $true, 1, $false, "29×92=", 2) | Out-Null
$d.Content.Find.Execute("83×89=", $true, $false, $false, $false, $false, $true, 1, $false, "27×15=", 2) | Out-Null
